# Apply updated LR-pair stats (natmi) per Dr Hou's revised normalization.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Sending cluster, Target cluster, then columns E..T in sheet order
$rowsData = @(
    @("ECs", "ECs", 3, 1, 3.194648666666667, 9.583946000000001, 0.1104372036365826, 0.122064847035162, 2, 0.6666666666666666, 0.2569163333333334, 0.770749, 0.2154602692984302, 0.2443643158836536, 0.8207574217282224, 7.386816795554001, 0.02379482963610366, 0.02982829283919019),
    @("ECs", "FAPs", 3, 1, 3.194648666666667, 9.583946000000001, 0.1104372036365826, 0.122064847035162, 3, 1, 0.5123676666666667, 1.537103, 0.4296919312505433, 0.4873352064520506, 1.636834683159778, 14.731512148438, 0.04745397531251268, 0.05948649743041867),
    @("ECs", "sCs", 3, 1, 3.194648666666667, 9.583946000000001, 0.1104372036365826, 0.122064847035162, 2, 1, 0.423123, 0.846246, 0.3548477994510263, 0.2683004776642958, 1.351729327786, 8.110375966716001, 0.0391883986879662, 0.03275005676555317),
    @("FAPs", "ECs", 3, 1, 16.97794966666667, 50.93384900000001, 0.5869181497900705, 0.64871322178746, 2, 0.6666666666666666, 0.2569163333333334, 0.770749, 0.2154602692984302, 0.2443643158836536, 4.36191257587789, 39.257213182901, 0.126457542609905, 0.1585223626467735),
    @("FAPs", "FAPs", 3, 1, 16.97794966666667, 50.93384900000001, 0.5869181497900705, 0.64871322178746, 3, 1, 0.5123676666666667, 1.537103, 0.4296919312505433, 0.4873352064520506, 8.698952455494114, 78.29057209944702, 0.252193993269291, 0.3161407918679667),
    @("FAPs", "sCs", 3, 1, 16.97794966666667, 50.93384900000001, 0.5869181497900705, 0.64871322178746, 2, 1, 0.423123, 0.846246, 0.3548477994510263, 0.2683004776642958, 7.183760996809002, 43.10256598085401, 0.2082666139108744, 0.1740500672727198),
    @("M1", "ECs", 3, 1, 0.1410116666666667, 0.423035, 0.004874693830745885, 0.005387937553646458, 2, 0.6666666666666666, 0.2569163333333334, 0.770749, 0.2154602692984302, 0.2443643158836536, 0.03622820035722223, 0.326053803215, 0.001050302845519905, 0.001316619674320663),
    @("M1", "FAPs", 3, 1, 0.1410116666666667, 0.423035, 0.004874693830745885, 0.005387937553646458, 3, 1, 0.5123676666666667, 1.537103, 0.4296919312505433, 0.4873352064520506, 0.07224981862277778, 0.650248367605, 0.002094616606388308, 0.002625731660057053),
    @("M1", "sCs", 3, 1, 0.1410116666666667, 0.423035, 0.004874693830745885, 0.005387937553646458, 2, 1, 0.423123, 0.846246, 0.3548477994510263, 0.2683004776642958, 0.05966527943500001, 0.35799167661, 0.001729774378837671, 0.001445586219268742),
    @("M2", "ECs", 3, 1, 0.3470173333333333, 1.041052, 0.01199619360545975, 0.01325924135378574, 2, 0.6666666666666666, 0.2569163333333334, 0.770749, 0.2154602692984302, 0.2443643158836536, 0.08915442088311111, 0.8023897879479999, 0.002584703104788464, 0.003240085442554102),
    @("M2", "FAPs", 3, 1, 0.3470173333333333, 1.041052, 0.01199619360545975, 0.01325924135378574, 3, 1, 0.5123676666666667, 1.537103, 0.4296919312505433, 0.4873352064520506, 0.1778004613728889, 1.600204152356, 0.005154667597985416, 0.006461695122544743),
    @("M2", "sCs", 3, 1, 0.3470173333333333, 1.041052, 0.01199619360545975, 0.01325924135378574, 2, 1, 0.423123, 0.846246, 0.3548477994510263, 0.2683004776642958, 0.146831015132, 0.880986090792, 0.004256822902685865, 0.003557460788686898),
    @("sCs", "ECs", 2, 1, 8.2666595, 16.533319, 0.2857737591371413, 0.2105747522699458, 2, 0.6666666666666666, 0.2569163333333334, 0.770749, 0.2154602692984302, 0.2443643158836536, 2.123839847655167, 12.743039085931, 0.06157289110211319, 0.05145695528081512),
    @("sCs", "FAPs", 2, 1, 8.2666595, 16.533319, 0.2857737591371413, 0.2105747522699458, 3, 1, 0.5123676666666667, 1.537103, 0.4296919312505433, 0.4873352064520506, 4.235569039142833, 25.413414234857, 0.1227946784643658, 0.1026204903710634),
    @("sCs", "sCs", 2, 1, 8.2666595, 16.533319, 0.2857737591371413, 0.2105747522699458, 2, 1, 0.423123, 0.846246, 0.3548477994510263, 0.2683004776642958, 3.4978137676185, 13.991255070474, 0.1014061895706622, 0.0564973066180672)
)

$r = 2
foreach ($row in $rowsData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "Nlgn2"
    $ws.Cells.Item($r, 3).Value = "Nrxn2"
    $ws.Cells.Item($r, 4).Value = $row[1]
    for ($i = 2; $i -lt $row.Length; $i++) {
        $ws.Cells.Item($r, $i + 3).Value = $row[$i]
    }
    $r = $r + 1
}
